# Applies the scheduled-runner Sheets update described in the commit diff.
# For each affected (sheet, cell): set the new numeric value, or clear the
# cell entirely when the target workbook no longer has that cell at all.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2293.111
$ws.Range("I9").Value = 2874
$ws.Range("J9").Value = 260
$ws.Range("K9").Value = 2874
$ws.Range("L9").Value = 260
$ws.Range("M9").Value = -2705
$ws.Range("N9").Value = -598
$ws.Range("H12").Value = 299.15384
$ws.Range("I12").Value = 299.08334
$ws.Range("K12").Value = 299.08334
$ws.Range("M12").Value = -129.08334
$ws.Range("H18").Value = 349.5
$ws.Range("I18").Value = 349.5
$ws.Range("K18").Value = 349.5
$ws.Range("M18").Value = -65.5
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H29").Value = 2522.4
$ws.Range("I29").Value = 1903
$ws.Range("J29").Value = 5000
$ws.Range("K29").Value = 5709
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = -5428
$ws.Range("N29").Value = -15562
$ws.Range("H42").Value = 1752
$ws.Range("J42").Value = 1000
$ws.Range("L42").Value = 3000
$ws.Range("N42").Value = -3460
$ws.Range("H50").Value = 381.25
$ws.Range("J50").Value = 381.25
$ws.Range("L50").Value = 1143.75
$ws.Range("N50").Value = -2093.75
$ws.Range("H58").Value = 2220.6667
$ws.Range("J58").Value = 6474.25
$ws.Range("L58").Value = 19422.75
$ws.Range("N58").Value = -19722.75
$ws.Range("H99").Value = 90913144
$ws.Range("J99").Value = 333345820
$ws.Range("L99").Value = 1000037460
$ws.Range("N99").Value = -1000040456
$ws.Range("H100").Value = 2000
$ws.Range("I100").Value = 1568.1818
$ws.Range("J100").Value = 4375
$ws.Range("K100").Value = 1568.1818
$ws.Range("L100").Value = 4375
$ws.Range("M100").Value = -1027.1818
$ws.Range("N100").Value = -5457
$ws.Range("H101").Value = 3096.4
$ws.Range("I101").Value = 3096.4
$ws.Range("K101").Value = 9289.200000000001
$ws.Range("M101").Value = -7667.200000000001
$ws.Range("H103").Value = 249.81818
$ws.Range("I103").Value = 224.8
$ws.Range("J103").Value = 500
$ws.Range("K103").Value = 674.4000000000001
$ws.Range("L103").Value = 1500
$ws.Range("M103").Value = -88.40000000000009
$ws.Range("N103").Value = -2672
$ws.Range("H105").Value = 50610
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H106").Value = 1287.7273
$ws.Range("I106").Value = 1287.7273
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 1287.7273
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -656.7273
$ws.Range("N106").ClearContents()
$ws.Range("H107").Value = 1358.0555
$ws.Range("I107").Value = 1358.0555
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1358.0555
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 561.9445000000001
$ws.Range("N107").ClearContents()
$ws.Range("H111").Value = 5596.2
$ws.Range("I111").Value = 5245.25
$ws.Range("K111").Value = 15735.75
$ws.Range("M111").Value = -12668.75
$ws.Range("H112").Value = 2692.96
$ws.Range("J112").Value = 2848.913
$ws.Range("L112").Value = 8546.739
$ws.Range("N112").Value = -10762.739
$ws.Range("H115").Value = 555
$ws.Range("I115").Value = 555
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 1665
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -98
$ws.Range("N115").ClearContents()
$ws.Range("H118").Value = 223.5
$ws.Range("I118").Value = 223.5
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 670.5
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 986.5
$ws.Range("N118").ClearContents()
$ws.Range("H125").Value = 6532.857
$ws.Range("I125").Value = 6532.857
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 58795.713
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -56335.713
$ws.Range("N125").ClearContents()
$ws.Range("H127").Value = 12778
$ws.Range("I127").Value = 10471.583
$ws.Range("J127").Value = 26616.5
$ws.Range("K127").Value = 31414.749
$ws.Range("L127").Value = 79849.5
$ws.Range("M127").Value = -26454.749
$ws.Range("N127").Value = -89769.5
$ws.Range("H137").Value = 8253.6
$ws.Range("I137").Value = 5638.8945
$ws.Range("J137").Value = 12769.909
$ws.Range("K137").Value = 16916.6835
$ws.Range("L137").Value = 38309.727
$ws.Range("M137").Value = -14366.6835
$ws.Range("N137").Value = -43409.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3804.6667
$ws.Range("I2").Value = 1268.4166
$ws.Range("K2").Value = 1268.4166
$ws.Range("M2").Value = -1155.4166
$ws.Range("H11").Value = 33336834
$ws.Range("I11").Value = 50000000
$ws.Range("K11").Value = 50000000
$ws.Range("M11").Value = -49999856
$ws.Range("H32").Value = 13007385
$ws.Range("I32").Value = 16133508
$ws.Range("J32").Value = 6950521.5
$ws.Range("K32").Value = 16133508
$ws.Range("L32").Value = 6950521.5
$ws.Range("M32").Value = -16133221
$ws.Range("N32").Value = -6951095.5
$ws.Range("H61").Value = 2843.4146
$ws.Range("I61").Value = 2887
$ws.Range("K61").Value = 2887
$ws.Range("M61").Value = -2675
$ws.Range("H92").Value = 61666.332
$ws.Range("J92").Value = 61666.332
$ws.Range("L92").Value = 61666.332
$ws.Range("N92").Value = -66658.33199999999
$ws.Range("H116").Value = 3804.6667
$ws.Range("I116").Value = 1268.4166
$ws.Range("K116").Value = 1268.4166
$ws.Range("M116").Value = 1025.5834
$ws.Range("H122").Value = 2513.742
$ws.Range("I122").Value = 1927.96
$ws.Range("K122").Value = 5783.88
$ws.Range("M122").Value = -3333.88
$ws.Range("H136").Value = 2843.4146
$ws.Range("I136").Value = 2887
$ws.Range("K136").Value = 8661
$ws.Range("M136").Value = -6111
$ws.Range("H139").Value = 85500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3804.6667
$ws.Range("I3").Value = 1268.4166
$ws.Range("K3").Value = 1268.4166
$ws.Range("M3").Value = -1154.4166
$ws.Range("H107").Value = 3575884.2
$ws.Range("I107").Value = 4171236
$ws.Range("J107").Value = 3773
$ws.Range("K107").Value = 4171236
$ws.Range("L107").Value = 3773
$ws.Range("M107").Value = -4169316
$ws.Range("N107").Value = -7613

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 71446696
$ws.Range("I134").Value = 166677630
$ws.Range("K134").Value = 500032890
$ws.Range("M134").Value = -500030355

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 96851.14
$ws.Range("J37").Value = 96851.14
$ws.Range("L37").Value = 290553.42
$ws.Range("N37").Value = -290777.42

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 67.85714
$ws.Range("I2").Value = 62.5
$ws.Range("K2").Value = 62.5
$ws.Range("M2").Value = 50.5
$ws.Range("H126").Value = 15156561
$ws.Range("I126").Value = 22730516
$ws.Range("K126").Value = 68191548
$ws.Range("M126").Value = -68189078
$ws.Range("H132").Value = 41671564
$ws.Range("I132").Value = 66671504
$ws.Range("K132").Value = 200014512
$ws.Range("M132").Value = -200011982

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2825.0625
$ws.Range("I93").Value = 2354.7896
$ws.Range("K93").Value = 2354.7896
$ws.Range("M93").Value = -1106.7896

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 13898744
$ws.Range("I136").Value = 14713259
$ws.Range("K136").Value = 44139777
$ws.Range("M136").Value = -44137227
